$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305558323860168
$ws.Range("B1").Value = 3.724775791168213
$ws.Range("C1").Value = 4.044373035430908
$ws.Range("D1").Value = 2.85629415512085
$ws.Range("E1").Value = 1.051178216934204
